$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 3548.4
$ws.Range("I20").Value = 2435.5
$ws.Range("J20").Value = 8000
$ws.Range("K20").Value = 2435.5
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = -2205.5
$ws.Range("N20").Value = -8460
# Row 21
$ws.Range("H21").Value = 45629.75
$ws.Range("I21").Value = 55012.668
$ws.Range("J21").Value = 40000
$ws.Range("K21").Value = 55012.668
$ws.Range("L21").Value = 40000
$ws.Range("M21").Value = -54544.668
$ws.Range("N21").Value = -40936
# Row 23
$ws.Range("H23").Value = 45629.75
$ws.Range("I23").Value = 55012.668
$ws.Range("J23").Value = 40000
$ws.Range("K23").Value = 55012.668
$ws.Range("L23").Value = 40000
$ws.Range("M23").Value = -54778.668
$ws.Range("N23").Value = -40468
# Row 29
$ws.Range("H29").Value = 434.875
$ws.Range("I29").Value = 137.25
$ws.Range("J29").Value = 732.5
$ws.Range("K29").Value = 411.75
$ws.Range("L29").Value = 2197.5
$ws.Range("M29").Value = -130.75
$ws.Range("N29").Value = -2759.5
# Row 35
$ws.Range("H35").Value = 3548.4
$ws.Range("I35").Value = 2435.5
$ws.Range("J35").Value = 8000
$ws.Range("K35").Value = 2435.5
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = -2056.5
$ws.Range("N35").Value = -8758
# Row 58
$ws.Range("H58").Value = 959.4
$ws.Range("I58").Value = 299.14285
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 897.4285500000001
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -747.4285500000001
$ws.Range("N58").Value = -7800
# Row 82
$ws.Range("H82").Value = 7264.2856
$ws.Range("I82").Value = 425
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 1275
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -869
$ws.Range("N82").Value = -30812
# Row 85
$ws.Range("H85").Value = 7264.2856
$ws.Range("I85").Value = 425
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 1275
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = 129
$ws.Range("N85").Value = -32808
# Row 87
$ws.Range("H87").Value = 15905.5
$ws.Range("J87").Value = 15905.5
$ws.Range("L87").Value = 15905.5
$ws.Range("N87").Value = -18401.5
# Row 90
$ws.Range("H90").Value = 15905.5
$ws.Range("J90").Value = 15905.5
$ws.Range("L90").Value = 47716.5
$ws.Range("N90").Value = -60196.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 39478.555
$ws.Range("I2").Value = 44324.25
$ws.Range("J2").Value = 713
$ws.Range("K2").Value = 44324.25
$ws.Range("L2").Value = 713
$ws.Range("M2").Value = -44211.25
$ws.Range("N2").Value = -939
# Row 60
$ws.Range("H60").Value = 43500
$ws.Range("I60").Value = 43500
$ws.Range("K60").Value = 43500
$ws.Range("M60").Value = -42767
# Row 116
$ws.Range("H116").Value = 39478.555
$ws.Range("I116").Value = 44324.25
$ws.Range("J116").Value = 713
$ws.Range("K116").Value = 44324.25
$ws.Range("L116").Value = 713
$ws.Range("M116").Value = -42030.25
$ws.Range("N116").Value = -5301

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 39478.555
$ws.Range("I3").Value = 44324.25
$ws.Range("J3").Value = 713
$ws.Range("K3").Value = 44324.25
$ws.Range("L3").Value = 713
$ws.Range("M3").Value = -44210.25
$ws.Range("N3").Value = -941

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 13014.125
$ws.Range("J50").Value = 13014.125
$ws.Range("L50").Value = 13014.125
$ws.Range("N50").Value = -14264.125
# Row 51
$ws.Range("H51").Value = 26300.273
$ws.Range("J51").Value = 9922.556
$ws.Range("L51").Value = 9922.556
$ws.Range("N51").Value = -11394.556
# Row 58
$ws.Range("H58").Value = 1575353.1
$ws.Range("I58").Value = 13224.4375
$ws.Range("J58").Value = 3497973
$ws.Range("K58").Value = 13224.4375
$ws.Range("L58").Value = 3497973
$ws.Range("M58").Value = -13021.4375
$ws.Range("N58").Value = -3498379
# Row 59
$ws.Range("H59").Value = 17748.25
$ws.Range("J59").Value = 17748.25
$ws.Range("L59").Value = 17748.25
$ws.Range("N59").Value = -20038.25
# Row 60
$ws.Range("H60").Value = 15740.134
$ws.Range("J60").Value = 9610.200000000001
$ws.Range("L60").Value = 9610.200000000001
$ws.Range("N60").Value = -10632.2
# Row 61
$ws.Range("H61").Value = 26300.273
$ws.Range("J61").Value = 9922.556
$ws.Range("L61").Value = 9922.556
$ws.Range("N61").Value = -10618.556
# Row 68
$ws.Range("H68").Value = 18110.334
$ws.Range("J68").Value = 18110.334
$ws.Range("L68").Value = 18110.334
$ws.Range("N68").Value = -19608.334
# Row 71
$ws.Range("H71").Value = 18110.334
$ws.Range("J71").Value = 18110.334
$ws.Range("L71").Value = 54331.00199999999
$ws.Range("N71").Value = -61819.00199999999
# Row 74
$ws.Range("H74").Value = 17973.334
$ws.Range("J74").Value = 19490.455
$ws.Range("L74").Value = 19490.455
$ws.Range("N74").Value = -21238.455
# Row 77
$ws.Range("H77").Value = 17973.334
$ws.Range("J77").Value = 19490.455
$ws.Range("L77").Value = 58471.36500000001
$ws.Range("N77").Value = -67207.36500000001
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 132
$ws.Range("H132").Value = 1425.8684
$ws.Range("I132").Value = 1134.5807
$ws.Range("J132").Value = 2715.8572
$ws.Range("K132").Value = 3403.7421
$ws.Range("L132").Value = 8147.571599999999
$ws.Range("M132").Value = -873.7420999999999
$ws.Range("N132").Value = -13207.5716
# Row 134
$ws.Range("H134").Value = 1382217.5
$ws.Range("I134").Value = 1940.9524
$ws.Range("J134").Value = 5005443.5
$ws.Range("K134").Value = 5822.857199999999
$ws.Range("L134").Value = 15016330.5
$ws.Range("M134").Value = -3287.857199999999
$ws.Range("N134").Value = -15021400.5
# Row 136
$ws.Range("H136").Value = 1575353.1
$ws.Range("I136").Value = 13224.4375
$ws.Range("J136").Value = 3497973
$ws.Range("K136").Value = 39673.3125
$ws.Range("L136").Value = 10493919
$ws.Range("M136").Value = -37123.3125
$ws.Range("N136").Value = -10499019

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 27
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5332

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 10725
$ws.Range("I69").Value = 9500
$ws.Range("J69").Value = 11950
$ws.Range("K69").Value = 9500
$ws.Range("L69").Value = 11950
$ws.Range("M69").Value = -8751
$ws.Range("N69").Value = -13448
# Row 72
$ws.Range("H72").Value = 10725
$ws.Range("I72").Value = 9500
$ws.Range("J72").Value = 11950
$ws.Range("K72").Value = 28500
$ws.Range("L72").Value = 35850
$ws.Range("M72").Value = -24756
$ws.Range("N72").Value = -43338
